# Update "Sprint Backlog Burndown.xlsx" to reflect the team's efforts for this sprint.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# --- Re-assign "Vitor" tasks to "Vitor/Jason" (rows 10-20, column B) ---
for ($r = 10; $r -le 20; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value() -eq "Vitor") {
        $cell.Value = "Vitor/Jason"
    }
}

# --- Update weekly remaining-work figures (columns D:G) to reflect progress ---
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 1

$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 2

$ws.Range("D5").Value = 2

$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1

$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 1

$ws.Range("D8").Value = 6
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 2

$ws.Range("D9").Value = 4
$ws.Range("G9").Value = 2

$ws.Range("D10").Value = 25
$ws.Range("E10").Value = 20
$ws.Range("F10").Value = 12
$ws.Range("G10").Value = 7

$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 2

$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 2

$ws.Range("D19").Value = 2

$ws.Range("D20").Value = 2

$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 1

$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1

# Recalculate totals / chart caches
$excel.Calculate()

# Move the active selection to the grand-total cell, like the author left it
$ws.Range("D31").Select()
